# "Generate Report for Archive"
#
# 1. The shared "Ready for handoff" status string becomes "In Translation"
#    everywhere it's used (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2).
# 2. The two "Latest Handoff/Handback Datetime"-ish status columns that used
#    to be sized for "Ready for handoff" (wide) shrink to fit the shorter
#    "In Translation" text: Overview columns E & F, and column C on the
#    zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns("E:F").ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns("C:C").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns("C:C").ColumnWidth = 12.5
